# Insert a new daily price record for "Feria Lagunitas de Puerto Montt -
# Frutilla" as row 34, pushing the existing rows 34-137 down to 35-138
# (matches the diff: dimension A1:T137 -> A1:T138, and every row from 34
# onward taking on the values that used to belong to the row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34..137 down by one, leaving a blank row 34 to populate.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Cells.Item(34, 1).Value  = 4
$ws.Cells.Item(34, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value  = "Los Lagos"
$ws.Cells.Item(34, 4).Value  = 44526
$ws.Cells.Item(34, 5).Value  = 10
$ws.Cells.Item(34, 6).Value  = "Fruta"
$ws.Cells.Item(34, 7).Value  = 100101
$ws.Cells.Item(34, 8).Value  = "Berries"
$ws.Cells.Item(34, 9).Value  = 100112025
$ws.Cells.Item(34, 10).Value = "Frutilla"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 1200
$ws.Cells.Item(34, 14).Value = 9500
$ws.Cells.Item(34, 15).Value = 9500
$ws.Cells.Item(34, 16).Value = 9500
$ws.Cells.Item(34, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(34, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(34, 19).Value = 1357
$ws.Cells.Item(34, 20).Value = 7
